# "Version 2." -> "Version 1." while preserving the spellcheck markers
# (w:proofErr spellStart/spellEnd) around "Version" and the _GoBack
# bookmark, and dropping the separate trailing "." run.
#
# Original runs / text layout (character offsets into $d.Content):
#   [0,5)  "Versi"   (run 1)
#   [5,7)  "on"      (run 2)
#   [7,9)  " 2"       (run 3, after spellEnd)
#   [9,10) "."        (run 4, after the _GoBack bookmark)
#
# We edit from right to left so earlier offsets stay valid, merging
# "Versi"+"on" into a single "Version" run and folding the "." into the
# " 2" -> " 1." run instead of leaving it as its own trailing run.

$d = $word.ActiveDocument

# Drop the standalone trailing "." run (its text becomes empty and Word
# removes the now-empty run).
$dot = $d.Range(9, 10)
$dot.Text = ""

# " 2" -> " 1." (keeps this run, now carries the final period too).
$verNum = $d.Range(7, 9)
$verNum.Text = " 1."

# Remove the second "Version" run ("on") so only "Versi" remains to be
# expanded below; this keeps the spellStart/spellEnd pair around a single
# run instead of two.
$on = $d.Range(5, 7)
$on.Text = ""

# "Versi" -> "Version"
$versi = $d.Range(0, 5)
$versi.Text = "Version"
